$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the header units for the VL_* columns from mV to V
$ws.Range("F1").Value = "VL_-20mA (V)"
$ws.Range("G1").Value = "VL_-10mA (V)"
$ws.Range("H1").Value = "VL_10mA (V)"
$ws.Range("I1").Value = "VL_20mA (V)"
$ws.Range("J1").Value = "VL_err (V)"

# Convert the VL_err column values from mV (0.1) to V (0.001)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.001
}

# Restore the view state (scrolled/selected range) as saved by the author
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("J2:J16").Select()
